# Append a new data row (row 63) with the latest reading pulled from Adafruit IO,
# matching the existing sheet layout: Timestamp, Feed Key, Value, Latitude, Longitude, Elevation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 63

$ws.Cells.Item($newRow, 1).Value = "2024-09-25T18:06:40Z"
$ws.Cells.Item($newRow, 2).Value = "temperature"

# Value column holds numeric-looking text ("25") in the existing data, so force
# the cell to Text format before assigning it to keep it a string, not a number.
$ws.Cells.Item($newRow, 3).NumberFormat = "@"
$ws.Cells.Item($newRow, 3).Value = "25"

$ws.Cells.Item($newRow, 4).Value = "N/A"
$ws.Cells.Item($newRow, 5).Value = "N/A"
$ws.Cells.Item($newRow, 6).Value = "N/A"
